$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Make room for the 3 new worker rows (table currently spans rows 16-19) ---
$ws.Rows("20:22").Insert()

# Row 22 becomes the new "last" row of the table - copy its special bottom-border
# formatting from row 19 (which still carries the old "last row" style at this point).
$ws.Range("B19:J19").Copy()
$ws.Range("B22:J22").PasteSpecial(-4122)   # xlPasteFormats

# Rows 19, 20 and 21 become ordinary interior rows - copy the "normal" row formatting
# from row 18.
$ws.Range("B18:J18").Copy()
$ws.Range("B19:J19").PasteSpecial(-4122)
$ws.Range("B20:J20").PasteSpecial(-4122)
$ws.Range("B21:J21").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Header / summary values ---
$ws.Range("D2").Value = "ESTADO DE CUENTA"
$ws.Range("E11").Value = 401980
$ws.Range("C13").Value = 7
$ws.Range("F13").Value = 4
$ws.Range("H15").Value = "Novedad de Ingreso"
$ws.Range("I15").Value = "Novedad de Retiro"
$ws.Range("J15").Value = "Observaciones"

# --- Worker table (rows 16-22) ---
$ws.Range("B16").Value = "CC"
$ws.Range("C16").Value = "9298939"
$ws.Range("D16").Value = "JAMEL ENRIQUE PANTOJA PIÑA"
$ws.Range("E16").Value = "2505"
$ws.Range("F16").Value = 62634
$ws.Range("G16").Value = 1565850

$ws.Range("B17").Value = "CC"
$ws.Range("C17").Value = "3805709"
$ws.Range("D17").Value = "JORGE LUIS QUINTANA MARTINEZ"
$ws.Range("E17").Value = "2505"
$ws.Range("F17").Value = 68328
$ws.Range("G17").Value = 1050901

$ws.Range("B18").Value = "CC"
$ws.Range("C18").Value = "9314008"
$ws.Range("D18").Value = "WALTER DE LA CRUZ ASENCIO CHAMORRO"
$ws.Range("E18").Value = "2505"
$ws.Range("F18").Value = 56940
$ws.Range("G18").Value = 1423500

$ws.Range("B19").Value = "CC"
$ws.Range("C19").Value = "1044918469"
$ws.Range("D19").Value = "CRISTOBAL JOSE PACHECO CORREA"
$ws.Range("E19").Value = "2505"
$ws.Range("F19").Value = 62634
$ws.Range("G19").Value = 1565850

$ws.Range("B20").Value = "CC"
$ws.Range("C20").Value = "11165994"
$ws.Range("D20").Value = "ALEXANDER REALES RAMIREZ"
$ws.Range("E20").Value = "2003"
$ws.Range("F20").Value = 59348
$ws.Range("G20").Value = 1483712

$ws.Range("B21").Value = "CC"
$ws.Range("C21").Value = "26162681"
$ws.Range("D21").Value = "MARGENIA DE JESUS BURGOS PEREZ"
$ws.Range("E21").Value = "1911"
$ws.Range("F21").Value = 31249
$ws.Range("G21").Value = 828116

$ws.Range("B22").Value = "CC"
$ws.Range("C22").Value = "1002188901"
$ws.Range("D22").Value = "JESICA GRACIELA MARTINEZ DIAZ"
$ws.Range("E22").Value = "2503"
$ws.Range("F22").Value = 60847
$ws.Range("G22").Value = 1521178
